# Append a new data row (row 91) to the "AYKO" sheet, mirroring the
# formatting of the existing data rows (text cells for most columns,
# numeric cells for Attachments / Coordenada_X / Coordenada_Y).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 91

# Force the text-like columns to Text format BEFORE assigning the
# values so Excel does not auto-convert numeric-looking strings
# (e.g. "6399", "9", "808258198") or date-looking strings
# (e.g. "7/14/2025") into numbers / dates.
# NOTE: use ${row} (braces) rather than $row immediately followed by
# ":" in a double-quoted string -- "$row:H" would otherwise be parsed
# as a PowerShell scoped-variable reference ($row:H) instead of the
# literal text "91:H".
$ws.Range("A${row}:H${row}").NumberFormat = "@"
$ws.Range("J${row}:L${row}").NumberFormat = "@"
$ws.Range("O${row}:P${row}").NumberFormat = "@"

$ws.Range("A${row}").Value = "6399"
$ws.Range("B${row}").Value = "7/14/2025"
$ws.Range("C${row}").Value = "ESCALADA AV. 966"
$ws.Range("D${row}").Value = "9"
$ws.Range("E${row}").Value = "808258198"
$ws.Range("F${row}").Value = "AYKO"
$ws.Range("G${row}").Value = "Pendiente"
$ws.Range("H${row}").Value = "Picada"
$ws.Range("I${row}").Value = 1
$ws.Range("J${row}").Value = "Cambio"
$ws.Range("K${row}").Value = "Sin equipos"
$ws.Range("L${row}").Value = "Pasante"
$ws.Range("M${row}").Value = -58.493069
$ws.Range("N${row}").Value = -34.646557
$ws.Range("O${row}").Value = "Devoto"
$ws.Range("P${row}").Value = "Capital Norte"
